# Add withdrawal report rows to "money_transfers" to get correct buy-dates
# of forex for potential transactions in the future, and update the two
# existing currency_conversions amounts used as the source data for this
# example.

$wb = $excel.ActiveWorkbook

# --- currency_conversions (sheet6) -----------------------------------
$cc = $wb.Worksheets.Item("currency_conversions")
$cc.Range("B2").Value = 4000
$cc.Range("B3").Value = 1000

# --- money_transfers (sheet7) -----------------------------------------
$mt = $wb.Worksheets.Item("money_transfers")

# row 3: this transfer is now a USD withdrawal of 4000, no buy-date yet
$mt.Range("B3").ClearContents()
$mt.Range("C3").Value = -4000
$mt.Range("D3").Value = 2
$mt.Range("E3").Value = "USD"

# row 4: withdrawal of 2000 USD with 2 fees, no buy-date yet
$mt.Range("B4").ClearContents()
$mt.Range("C4").Value = -2000
$mt.Range("D4").Value = 2
$mt.Range("E4").Value = "USD"

# row 5: withdrawal of 2000 EUR, no buy-date yet
$mt.Range("B5").ClearContents()
$mt.Range("C5").Value = -2000
$mt.Range("E5").Value = "EUR"

# new column G placeholders on rows 4 and 5 (blank cells, formatted like F2)
$mt.Range("G4").Value = 0
$mt.Range("G4").Font.Color = 0
$mt.Range("G4").ClearContents()

$mt.Range("G5").Value = 0
$mt.Range("G5").Font.Color = 0
$mt.Range("G5").ClearContents()
